# Update "想去人数" (want-to-go count) figures on each sheet to match the
# newly generated output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 13023
$ws.Range("F8").Value = 55
$ws.Range("F10").Value = 260
$ws.Range("F11").Value = 3062
$ws.Range("F13").Value = 6475
$ws.Range("F16").Value = 3414
$ws.Range("F18").Value = 163
$ws.Range("F20").Value = 39
$ws.Range("F21").Value = 64
$ws.Range("F23").Value = 42
$ws.Range("F24").Value = 3612
$ws.Range("F27").Value = 2781
$ws.Range("F28").Value = 2781
$ws.Range("F29").Value = 412
$ws.Range("F30").Value = 1886
$ws.Range("F33").Value = 6637
$ws.Range("F36").Value = 847
$ws.Range("F37").Value = 1983
$ws.Range("F40").Value = 1040
$ws.Range("F45").Value = 1141
$ws.Range("F46").Value = 135
$ws.Range("F47").Value = 1201
$ws.Range("F48").Value = 1784
$ws.Range("F49").Value = 160
$ws.Range("F50").Value = 1172

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 125

# --- 本地生活 (Local Life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 432
$ws.Range("F3").Value = 598

# --- 全部类型 (All Types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 432
$ws.Range("F7").Value = 598
$ws.Range("F9").Value = 13023
$ws.Range("F12").Value = 260
$ws.Range("F13").Value = 3062
$ws.Range("F15").Value = 6475
$ws.Range("F17").Value = 3414
$ws.Range("F19").Value = 163
$ws.Range("F21").Value = 39
$ws.Range("F22").Value = 64
$ws.Range("F25").Value = 42
$ws.Range("F26").Value = 3612
$ws.Range("F28").Value = 2781
$ws.Range("F29").Value = 412
$ws.Range("F30").Value = 1886
$ws.Range("F33").Value = 6637
$ws.Range("F37").Value = 847
$ws.Range("F38").Value = 1983
$ws.Range("F42").Value = 1040
$ws.Range("F46").Value = 135
$ws.Range("F48").Value = 1784
$ws.Range("F50").Value = 160
